# Apply the "window" step update to the worksheet:
# - Add a new header "window" in E1
# - Add a new data row (row 10): Lp.=9, Nazwa=GOOGLE, Strona=https://www.google.com/
# - Update the active selection to B10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("E1").Value = "window"

# New data row
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "GOOGLE"
$ws.Cells.Item(10, 3).Value = "https://www.google.com/"

# Update selection to match the recorded sheet view state
$ws.Range("B10").Select()
